# Actualización automática 2025-07-04 17:20:07
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("E5").Value = 401.81
$wsVentasPorGrupo.Range("M5").Value = 3969.29
$wsVentasPorGrupo.Range("E22").Value = "1 de 20"
$wsVentasPorGrupo.Range("M22").Value = "5 de 20"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F5").Value = 4371.1
$wsVentaMensual.Range("F22").Value = 19235.06

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D4").Value = 401.81
$wsCumplimiento.Range("E4").Value = 112.021046659336
$wsCumplimiento.Range("F4").Value = 0.7819885594931661

$wsCumplimiento.Range("D16").Value = 15879.11
$wsCumplimiento.Range("E16").Value = 28387.13
$wsCumplimiento.Range("F16").Value = 0.3587182918630542

$wsCumplimiento.Range("D19").Value = 19235.06
$wsCumplimiento.Range("E19").Value = 46142.93762291769
$wsCumplimiento.Range("F19").Value = 0.2942130487223322
